$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.837.54"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").Value = "2.657.64"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.631"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.127"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.398"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.87%  "

$ws.Range("E11").Value = "  -0.88%  "

$ws.Range("E12").Value = "  +1.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.21%  "

$ws.Range("E14").Value = "  +0.48%  "

$ws.Range("D15").Value = "3.134.19"
$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("D16").Value = "65.654.65"
$ws.Range("E16").Value = "  +0.36%  "

$ws.Range("D17").Value = "2.644.13"
$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("E18").Value = "  -0.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.30%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000113"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.52%  "

$ws.Range("E27").Value = "  +1.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "554.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.52%  "

$ws.Range("E29").Value = "  -2.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("E31").Value = "  -1.60%  "

$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("E33").Value = "  +1.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("E35").Value = "  -0.94%  "

$ws.Range("E36").Value = "  -0.83%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("E39").Value = "  +0.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "161.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.91%  "

$ws.Range("E45").Value = "  -0.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.11%  "

$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("E48").Value = "  -0.70%  "

$ws.Range("E49").Value = "  -1.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.70%  "

$ws.Range("E51").Value = "  +7.16%  "
